# Scheduled-runner refresh: update cached market-board pricing / profit
# figures (columns H-N) across several leve-crafting sheets. No formulas
# are involved -- every touched cell holds a static number pasted in by
# the runner, so we just overwrite the stale values with the refreshed
# ones (and add/remove a couple of profit cells where the leve moved
# in/out of profitability).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 2966.6667
$ws.Cells.Item(74, 9).Value = 2966.6667
$ws.Cells.Item(74, 11).Value = 2966.6667
$ws.Cells.Item(74, 13).Value = -2030.6667

$ws.Cells.Item(77, 8).Value = 2966.6667
$ws.Cells.Item(77, 9).Value = 2966.6667
$ws.Cells.Item(77, 11).Value = 14833.3335
$ws.Cells.Item(77, 13).Value = -10153.3335

$ws.Cells.Item(80, 8).Value = 472.5
$ws.Cells.Item(80, 9).Value = 247.5
$ws.Cells.Item(80, 10).Value = 697.5
$ws.Cells.Item(80, 11).Value = 742.5
$ws.Cells.Item(80, 12).Value = 2092.5
$ws.Cells.Item(80, 13).Value = 255.5
$ws.Cells.Item(80, 14).Value = -4088.5

$ws.Cells.Item(83, 8).Value = 472.5
$ws.Cells.Item(83, 9).Value = 247.5
$ws.Cells.Item(83, 10).Value = 697.5
$ws.Cells.Item(83, 11).Value = 2227.5
$ws.Cells.Item(83, 12).Value = 6277.5
$ws.Cells.Item(83, 13).Value = 2764.5
$ws.Cells.Item(83, 14).Value = -16261.5

$ws.Cells.Item(132, 8).Value = 1070.4445
$ws.Cells.Item(132, 9).Value = 1070.4445
$ws.Cells.Item(132, 11).Value = 3211.3335
$ws.Cells.Item(132, 13).Value = -681.3335000000002

$ws.Cells.Item(137, 8).Value = 1754.7941
$ws.Cells.Item(137, 9).Value = 1362.826
$ws.Cells.Item(137, 11).Value = 4088.478
$ws.Cells.Item(137, 13).Value = -1538.478

$ws.Cells.Item(138, 8).Value = 3494.5083
$ws.Cells.Item(138, 10).Value = 3544.9636
$ws.Cells.Item(138, 12).Value = 10634.8908
$ws.Cells.Item(138, 14).Value = -20914.8908

$ws.Cells.Item(141, 8).Value = 6523.125
$ws.Cells.Item(141, 9).Value = 6259.385
$ws.Cells.Item(141, 11).Value = 18778.155
$ws.Cells.Item(141, 13).Value = -13598.155

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1286.5714
$ws.Cells.Item(2, 10).Value = 1450
$ws.Cells.Item(2, 12).Value = 1450
$ws.Cells.Item(2, 14).Value = -1676

$ws.Cells.Item(32, 8).Value = 5086.175
$ws.Cells.Item(32, 9).Value = 3800.5
$ws.Cells.Item(32, 11).Value = 3800.5
$ws.Cells.Item(32, 13).Value = -3513.5

$ws.Cells.Item(44, 8).Value = 39996.668
$ws.Cells.Item(44, 10).Value = 39996.668
$ws.Cells.Item(44, 12).Value = 39996.668
$ws.Cells.Item(44, 14).Value = -40972.668

$ws.Cells.Item(55, 8).Value = 39996.668
$ws.Cells.Item(55, 10).Value = 39996.668
$ws.Cells.Item(55, 12).Value = 39996.668
$ws.Cells.Item(55, 14).Value = -40626.668

$ws.Cells.Item(74, 8).Value = 1364.0714
$ws.Cells.Item(74, 9).Value = 1016.4167
$ws.Cells.Item(74, 11).Value = 1016.4167
$ws.Cells.Item(74, 13).Value = -142.4167

$ws.Cells.Item(77, 8).Value = 1364.0714
$ws.Cells.Item(77, 9).Value = 1016.4167
$ws.Cells.Item(77, 11).Value = 5082.0835
$ws.Cells.Item(77, 13).Value = -714.0834999999997

$ws.Cells.Item(116, 8).Value = 1286.5714
$ws.Cells.Item(116, 10).Value = 1450
$ws.Cells.Item(116, 12).Value = 1450
$ws.Cells.Item(116, 14).Value = -6038

$ws.Cells.Item(122, 8).Value = 4064.5557
$ws.Cells.Item(122, 10).Value = 1473
$ws.Cells.Item(122, 12).Value = 4419
$ws.Cells.Item(122, 14).Value = -9319

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1286.5714
$ws.Cells.Item(3, 10).Value = 1450
$ws.Cells.Item(3, 12).Value = 1450
$ws.Cells.Item(3, 14).Value = -1678

$ws.Cells.Item(105, 8).Value = 6030.8
$ws.Cells.Item(105, 9).Value = 4994
$ws.Cells.Item(105, 11).Value = 4994
$ws.Cells.Item(105, 13).Value = -3247

$ws.Cells.Item(134, 8).Value = 2473.5881
$ws.Cells.Item(134, 9).Value = 2255.9167
$ws.Cells.Item(134, 10).Value = 2996
$ws.Cells.Item(134, 11).Value = 6767.750100000001
$ws.Cells.Item(134, 12).Value = 8988
$ws.Cells.Item(134, 13).Value = -4232.750100000001
$ws.Cells.Item(134, 14).Value = -14058

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1933
$ws.Cells.Item(31, 9).Value = 1625
$ws.Cells.Item(31, 11).Value = 1625
$ws.Cells.Item(31, 13).Value = -1330

$ws.Cells.Item(34, 8).Value = 1933
$ws.Cells.Item(34, 9).Value = 1625
$ws.Cells.Item(34, 11).Value = 1625
$ws.Cells.Item(34, 13).Value = -1423

$ws.Cells.Item(51, 8).Value = 49991.668
$ws.Cells.Item(51, 10).Value = 49991.668
$ws.Cells.Item(51, 12).Value = 49991.668
$ws.Cells.Item(51, 14).Value = -51463.668

$ws.Cells.Item(58, 8).Value = 2754.3333
$ws.Cells.Item(58, 9).Value = 2332.6667
$ws.Cells.Item(58, 10).Value = 2965.1667
$ws.Cells.Item(58, 11).Value = 2332.6667
$ws.Cells.Item(58, 12).Value = 2965.1667
$ws.Cells.Item(58, 13).Value = -2129.6667
$ws.Cells.Item(58, 14).Value = -3371.1667

$ws.Cells.Item(61, 8).Value = 49991.668
$ws.Cells.Item(61, 10).Value = 49991.668
$ws.Cells.Item(61, 12).Value = 49991.668
$ws.Cells.Item(61, 14).Value = -50687.668

$ws.Cells.Item(99, 8).Value = 7674.6665
$ws.Cells.Item(99, 9).Value = 9005
$ws.Cells.Item(99, 10).Value = 5014
$ws.Cells.Item(99, 11).Value = 9005
$ws.Cells.Item(99, 12).Value = 5014
$ws.Cells.Item(99, 13).Value = -7507
$ws.Cells.Item(99, 14).Value = -8010

$ws.Cells.Item(105, 8).Value = 2122.8
$ws.Cells.Item(105, 9).Value = 2223
$ws.Cells.Item(105, 11).Value = 2223
$ws.Cells.Item(105, 13).Value = -476

$ws.Cells.Item(107, 8).Value = 1931.45
$ws.Cells.Item(107, 9).Value = 1220.6666
$ws.Cells.Item(107, 10).Value = 2513
$ws.Cells.Item(107, 11).Value = 1220.6666
$ws.Cells.Item(107, 12).Value = 2513
$ws.Cells.Item(107, 13).Value = 699.3334
$ws.Cells.Item(107, 14).Value = -6353

$ws.Cells.Item(126, 8).Value = 7674.6665
$ws.Cells.Item(126, 9).Value = 9005
$ws.Cells.Item(126, 10).Value = 5014
$ws.Cells.Item(126, 11).Value = 27015
$ws.Cells.Item(126, 12).Value = 15042
$ws.Cells.Item(126, 13).Value = -24545
$ws.Cells.Item(126, 14).Value = -19982

$ws.Cells.Item(136, 8).Value = 2754.3333
$ws.Cells.Item(136, 9).Value = 2332.6667
$ws.Cells.Item(136, 10).Value = 2965.1667
$ws.Cells.Item(136, 11).Value = 6998.000100000001
$ws.Cells.Item(136, 12).Value = 8895.500100000001
$ws.Cells.Item(136, 13).Value = -4448.000100000001
$ws.Cells.Item(136, 14).Value = -13995.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 3000
$ws.Cells.Item(51, 9).Value = 3000
$ws.Cells.Item(51, 11).Value = 9000
$ws.Cells.Item(51, 13).Value = -8540

$ws.Cells.Item(138, 8).Value = 7646.7144
$ws.Cells.Item(138, 9).Value = 5882
$ws.Cells.Item(138, 11).Value = 17646
$ws.Cells.Item(138, 13).Value = -12506

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 7546.125
$ws.Cells.Item(80, 9).Value = 6795.4
$ws.Cells.Item(80, 11).Value = 6795.4
$ws.Cells.Item(80, 13).Value = -5797.4

$ws.Cells.Item(83, 8).Value = 7546.125
$ws.Cells.Item(83, 9).Value = 6795.4
$ws.Cells.Item(83, 11).Value = 33977
$ws.Cells.Item(83, 13).Value = -28985

$ws.Cells.Item(122, 8).Value = 1628.1666
$ws.Cells.Item(122, 9).Value = 1317.5
$ws.Cells.Item(122, 10).Value = 2249.5
$ws.Cells.Item(122, 11).Value = 3952.5
$ws.Cells.Item(122, 12).Value = 6748.5
$ws.Cells.Item(122, 13).Value = -1502.5
$ws.Cells.Item(122, 14).Value = -11648.5

$ws.Cells.Item(132, 8).Value = 4082.762
$ws.Cells.Item(132, 9).Value = 3885.6365
$ws.Cells.Item(132, 11).Value = 11656.9095
$ws.Cells.Item(132, 13).Value = -9126.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 13).ClearContents()
$ws.Cells.Item(55, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 4205.25
$ws.Cells.Item(132, 9).Value = 3607.3333
$ws.Cells.Item(132, 11).Value = 10821.9999
$ws.Cells.Item(132, 13).Value = -8291.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(125, 8).Value = 78775
$ws.Cells.Item(125, 10).Value = 78775
$ws.Cells.Item(125, 12).Value = 78775
$ws.Cells.Item(125, 14).Value = -88615
